# repull data, push all data, mean calculation
# Column F ("dSF") values were re-pulled for a subset of rows; update them
# in place to match the freshly-pulled source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Row=2;  Value=6}
    @{Row=3;  Value=-3}
    @{Row=5;  Value=-3}
    @{Row=8;  Value=-4}
    @{Row=9;  Value=-4}
    @{Row=12; Value=3}
    @{Row=13; Value=2}
    @{Row=16; Value=-2}
    @{Row=17; Value=-5}
    @{Row=18; Value=2}
    @{Row=20; Value=5}
    @{Row=21; Value=-1}
    @{Row=22; Value=-6}
    @{Row=23; Value=-2}
    @{Row=25; Value=4}
    @{Row=29; Value=3}
    @{Row=33; Value=0}
    @{Row=34; Value=-1}
    @{Row=35; Value=-1}
    @{Row=39; Value=-3}
    @{Row=42; Value=0}
    @{Row=44; Value=3}
    @{Row=49; Value=-1}
    @{Row=55; Value=0}
    @{Row=56; Value=-5}
    @{Row=57; Value=0}
    @{Row=58; Value=2}
    @{Row=59; Value=-3}
    @{Row=60; Value=-3}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 6).Value = $u.Value
}
